$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new senior-editor columns: Q (Caroline) and R (CandyC)
$ws.Range("Q1").Value = "Caroline"
$ws.Range("R1").Value = "CandyC"

# Copy header style (bold, bordered, centered) from P1 onto the new header cells
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the new Q/R columns with the default "10-19" shift for every data row
$ws.Range("Q2:R31").Value = "10-19"

# Update individual leave/shift entries in column H (Daisy)
$ws.Range("H2").Value = "half off"
$ws.Range("H3").Value = "13-22"
$ws.Range("H4").Value = "AL"
$ws.Range("H6").Value = "13-22"
$ws.Range("H7").Value = "13-22"
$ws.Range("H11").Value = "13-22"
$ws.Range("H12").Value = "13-22"
$ws.Range("H15").Value = "13-22"
$ws.Range("H16").Value = "13-22"
$ws.Range("H20").Value = "13-22"
$ws.Range("H21").Value = "13-22"
$ws.Range("H24").Value = "13-22"
$ws.Range("H25").Value = "13-22"
$ws.Range("H26").Value = "13-22"
$ws.Range("H29").Value = "13-22"
$ws.Range("H30").Value = "13-22"
